$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.059.66"
$ws.Range("E2").Value = "  -3.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.575.89"
$ws.Range("E3").Value = "  -2.53%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.44"
$ws.Range("E5").Value = "  -0.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.17"

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  +2.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.79"
$ws.Range("E9").Value = "  +3.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0996"
$ws.Range("E10").Value = "  -3.37%  "

$ws.Range("E12").Value = "  -2.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.033.07"
$ws.Range("E13").Value = "  -2.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "57.998.98"
$ws.Range("E14").Value = "  -3.11%  "

$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.567.25"
$ws.Range("E16").Value = "  -3.19%  "

$ws.Range("E17").Value = "  -2.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.39"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "332.98"
$ws.Range("E19").Value = "  -3.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.01"
$ws.Range("E20").Value = "  -2.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.16"
$ws.Range("E21").Value = "  -3.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.79"
$ws.Range("E23").Value = "  -0.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.419"
$ws.Range("E24").Value = "  +1.39%  "

$ws.Range("E25").Value = "  +0.55%  "

$ws.Range("E26").Value = "  -4.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.02"
$ws.Range("E27").Value = "  -3.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0723"
$ws.Range("E29").Value = "  -3.77%  "

$ws.Range("E30").Value = "  -1.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.77"
$ws.Range("E31").Value = "  +2.69%  "

$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("E34").Value = "  -3.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "36.85"
$ws.Range("E35").Value = "  -1.47%  "

$ws.Range("E36").Value = "  -4.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.831"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.817"
$ws.Range("E38").Value = "  -2.92%  "

$ws.Range("E39").Value = "  -4.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.57"
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "280.74"
$ws.Range("E41").Value = "  -3.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.64"
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("E44").Value = "  -3.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0946"
$ws.Range("E45").Value = "  -0.63%  "

$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.80"
$ws.Range("E47").Value = "  -1.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0226"
$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.908.62"
$ws.Range("E49").Value = "  -3.37%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.74"
$ws.Range("E50").Value = "  -4.90%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.39"
$ws.Range("E51").Value = "  -3.70%  "
